# Apply cryptos list update (GitHub Actions scrape refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.652.41"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "2.199.88"
$ws.Range("E3").Value = "  -1.56%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.25"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.615"
$ws.Range("E6").Value = "  -1.48%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.29"
$ws.Range("E7").Value = "  -2.71%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.400"
$ws.Range("E9").Value = "  -0.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "56.92"
$ws.Range("E10").Value = "  -3.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0888"
$ws.Range("E11").Value = "  +0.87%  "
$ws.Range("E12").Value = "  -0.87%  "
$ws.Range("D13").Value = "2.532.07"
$ws.Range("E13").Value = "  -1.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.33"
$ws.Range("E14").Value = "  -2.68%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.00"
$ws.Range("E15").Value = "  -0.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.792"
$ws.Range("E16").Value = "  -1.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.55"
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("D18").Value = "2.202.23"
$ws.Range("E18").Value = "  -0.96%  "
$ws.Range("D19").Value = "41.643.08"
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.81"
$ws.Range("E21").Value = "  -2.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.02"
$ws.Range("E22").Value = "  -0.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.96"
$ws.Range("E23").Value = "  -2.96%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.998"
$ws.Range("E24").Value = "  -0.29%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.36"
$ws.Range("E25").Value = "  +1.42%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.34"
$ws.Range("E26").Value = "  -1.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.57"
$ws.Range("E27").Value = "  +0.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "168.12"
$ws.Range("E29").Value = "  -4.19%  "
$ws.Range("E30").Value = "  +0.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.66"
$ws.Range("E31").Value = "  -2.35%  "
$ws.Range("E32").Value = "  -7.79%  "
$ws.Range("E33").Value = "  -1.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.94"
$ws.Range("E34").Value = "  -1.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.57"
$ws.Range("E35").Value = "  -1.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0643"
$ws.Range("E36").Value = "  +2.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.54"
$ws.Range("E37").Value = "  -6.71%  "
$ws.Range("E38").Value = "  -6.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.32"
$ws.Range("E39").Value = "  -2.46%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.000241"
$ws.Range("E40").Value = "  +0.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0239"
$ws.Range("E42").Value = "  +0.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.66"
$ws.Range("E43").Value = "  +0.34%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.20"
$ws.Range("E44").Value = "  +0.40%  "
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0951"
$ws.Range("E45").Value = "  -0.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "96.37"
$ws.Range("E46").Value = "  -3.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.35"
$ws.Range("E47").Value = "  -12.06%  "
$ws.Range("D48").Value = "1.451.79"
$ws.Range("E48").Value = "  -2.34%  "
$ws.Range("E49").Value = "  -0.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.01"
$ws.Range("E50").Value = "  -5.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.06"
$ws.Range("E51").Value = "  -1.85%  "
